$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 totals: split the single shared {Total} placeholder into
# per-column {Total1}..{Total11} placeholders (H14 keeps the renamed text,
# J14:S14 each get their own unique placeholder text). ---
$ws.Range("H14").Value = "{Total1}"
$ws.Range("J14").Value = "{Total2}"
$ws.Range("K14").Value = "{Total3}"
$ws.Range("L14").Value = "{Total4}"
$ws.Range("M14").Value = "{Total5}"
$ws.Range("N14").Value = "{Total6}"
$ws.Range("O14").Value = "{Total7}"
$ws.Range("P14").Value = "{Total8}"
$ws.Range("Q14").Value = "{Total9}"
$ws.Range("R14").Value = "{Total10}"
$ws.Range("S14").Value = "{Total11}"

# --- Row 13 (the column header row) now wraps its header text. ---
$ws.Range("A13:S13").WrapText = $true

# --- Move the active selection in the frozen (bottom) pane to T8. ---
$ws.Range("T8").Select()
